# Apply the cell-level updates captured in the commit diff for cryptos.xlsx
# (coin price/volume/hour refresh run by the scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) hold numeric-looking text in this sheet.
# Force the Text number format before writing so COM keeps the exact string
# (e.g. '5.460', '0.1040', '16') instead of silently coercing to a Double
# and dropping trailing zeros / precision.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '245.82'
$ws.Range("G2").Value = '16'
# Row 3
$ws.Range("G3").Value = '16'
# Row 4
$ws.Range("D4").Value = '5.460'
$ws.Range("G4").Value = '16'
# Row 5
$ws.Range("D5").Value = '0.05617'
$ws.Range("G5").Value = '16'
# Row 6
$ws.Range("D6").Value = '6.471'
$ws.Range("G6").Value = '16'
# Row 7
$ws.Range("D7").Value = '0.8049'
$ws.Range("G7").Value = '16'
# Row 8
$ws.Range("D8").Value = '1.052'
$ws.Range("G8").Value = '16'
# Row 9
$ws.Range("D9").Value = '0.1427'
$ws.Range("G9").Value = '16'
# Row 10
$ws.Range("D10").Value = '0.07372'
$ws.Range("G10").Value = '16'
# Row 11
$ws.Range("D11").Value = '0.03181'
$ws.Range("G11").Value = '16'
# Row 12
$ws.Range("B12").Value = 'ProBitToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D12").Value = '0.1260'
$ws.Range("E12").Value = '11ProBitTokenPROBBestin24h'
$ws.Range("G12").Value = '16'
# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.02937'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("G13").Value = '16'
# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09263'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("G14").Value = '16'
# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001669'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("G15").Value = '16'
# Row 16
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '3.270'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("G16").Value = '16'
# Row 17
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '0.04719'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("G17").Value = '16'
# Row 18
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '0.0005826'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("G18").Value = '16'
# Row 19
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").Value = '0.006421'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("G19").Value = '16'
# Row 20
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").Value = '0.005060'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("G20").Value = '16'
# Row 21
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").Value = '0.001055'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("G21").Value = '16'
# Row 22
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = '0.0001504'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("G22").Value = '16'
# Row 23
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '3.986'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("G23").Value = '16'
# Row 24
$ws.Range("B24").Value = 'GateToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D24").Value = '3.383'
$ws.Range("E24").Value = '23GateTokenGT'
$ws.Range("G24").Value = '16'
# Row 25
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.125'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("G25").Value = '16'
# Row 26
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '0.3295'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("G26").Value = '16'
# Row 27
$ws.Range("D27").Value = '0.0002908'
$ws.Range("G27").Value = '16'
# Row 28
$ws.Range("G28").Value = '16'
# Row 29
$ws.Range("G29").Value = '16'
# Row 30
$ws.Range("G30").Value = '16'
# Row 31
$ws.Range("G31").Value = '16'
# Row 32
$ws.Range("G32").Value = '16'
# Row 33
$ws.Range("G33").Value = '16'
# Row 34
$ws.Range("G34").Value = '16'
# Row 35
$ws.Range("G35").Value = '16'
# Row 36
$ws.Range("G36").Value = '16'
# Row 37
$ws.Range("G37").Value = '16'
# Row 38
$ws.Range("G38").Value = '16'
# Row 39
$ws.Range("G39").Value = '16'
# Row 40
$ws.Range("D40").Value = '0.04158'
$ws.Range("G40").Value = '16'
# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '0.006908'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("G41").Value = '16'
# Row 42
$ws.Range("D42").Value = '0.003509'
$ws.Range("G42").Value = '16'
# Row 43
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").Value = '0.1040'
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("G43").Value = '16'
# Row 44
$ws.Range("D44").Value = '0.009306'
$ws.Range("G44").Value = '16'
# Row 45
$ws.Range("D45").Value = '0.00005665'
$ws.Range("G45").Value = '16'
# Row 46
$ws.Range("G46").Value = '16'
# Row 47
$ws.Range("D47").Value = '0.6819'
$ws.Range("G47").Value = '16'
# Row 48
$ws.Range("D48").Value = '0.01625'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("G48").Value = '16'
# Row 49
$ws.Range("G49").Value = '16'
# Row 50
$ws.Range("G50").Value = '16'
# Row 51
$ws.Range("G51").Value = '16'
